$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "89.177.07"
$ws.Range("E2").Value = "  -3.75%  "

$ws.Range("D3").Value = "3.155.51"
$ws.Range("E3").Value = "  -4.38%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'215.67"
$ws.Range("E5").Value = "  -1.89%  "

$ws.Range("D6").Value = "'633.42"
$ws.Range("E6").Value = "  -0.29%  "

$ws.Range("D7").Value = "'0.398"
$ws.Range("E7").Value = "  -4.38%  "

$ws.Range("D8").Value = "'0.738"
$ws.Range("E8").Value = "  +1.16%  "

$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.15%  "

$ws.Range("D10").Value = "3.154.90"
$ws.Range("E10").Value = "  -4.23%  "

$ws.Range("D11").Value = "'0.556"
$ws.Range("E11").Value = "  -6.56%  "

$ws.Range("E12").Value = "  -1.02%  "

$ws.Range("D13").Value = "'0.0000252"
$ws.Range("E13").Value = "  -4.89%  "

$ws.Range("E14").Value = "  -1.07%  "

$ws.Range("D15").Value = "88.936.44"
$ws.Range("E15").Value = "  -3.57%  "

$ws.Range("D16").Value = "3.737.56"
$ws.Range("E16").Value = "  -4.20%  "

$ws.Range("D17").Value = "'32.50"
$ws.Range("E17").Value = "  -5.72%  "

$ws.Range("D18").Value = "3.140.31"
$ws.Range("E18").Value = "  -3.43%  "

$ws.Range("D19").Value = "'0.0000232"
$ws.Range("E19").Value = "  +21.54%  "

$ws.Range("E20").Value = "  +1.77%  "

$ws.Range("D21").Value = "'13.35"
$ws.Range("E21").Value = "  -5.81%  "

$ws.Range("D22").Value = "'427.71"
$ws.Range("E22").Value = "  -2.87%  "

$ws.Range("D23").Value = "'8.41"
$ws.Range("E23").Value = "  -6.18%  "

$ws.Range("D24").Value = "'4.90"
$ws.Range("E24").Value = "  -8.03%  "

$ws.Range("D25").Value = "'5.43"
$ws.Range("E25").Value = "  -0.08%  "

$ws.Range("D26").Value = "'11.57"
$ws.Range("E26").Value = "  -6.85%  "

$ws.Range("D27").Value = "'81.36"
$ws.Range("E27").Value = "  +5.57%  "

$ws.Range("D28").Value = "3.318.61"
$ws.Range("E28").Value = "  -4.59%  "

$ws.Range("D29").Value = "'1.01"
$ws.Range("E29").Value = "  +0.82%  "

$ws.Range("D30").Value = "'0.160"
$ws.Range("E30").Value = "  -10.62%  "

$ws.Range("D31").Value = "'0.955"
$ws.Range("E31").Value = "  -4.23%  "

$ws.Range("D32").Value = "'4.06"
$ws.Range("E32").Value = "  +10.73%  "

$ws.Range("D33").Value = "'8.25"
$ws.Range("E33").Value = "  -6.78%  "

$ws.Range("D34").Value = "'512.19"
$ws.Range("E34").Value = "  -8.50%  "

$ws.Range("D35").Value = "'7.12"
$ws.Range("E35").Value = "  -1.47%  "

$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "'1.30"
$ws.Range("E36").Value = "  +1.04%  "

$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.140"
$ws.Range("E37").Value = "  +6.45%  "

$ws.Range("E38").Value = "  -5.11%  "

$ws.Range("D39").Value = "'21.99"
$ws.Range("E39").Value = "  -3.46%  "

$ws.Range("D40").Value = "'22.25"
$ws.Range("E40").Value = "  -1.05%  "

$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.31%  "

$ws.Range("E42").Value = "  -0.04%  "

$ws.Range("D43").Value = "'1.87"
$ws.Range("E43").Value = "  -7.03%  "

$ws.Range("D44").Value = "'0.365"
$ws.Range("E44").Value = "  -7.78%  "

$ws.Range("D45").Value = "'145.94"
$ws.Range("E45").Value = "  -2.77%  "

$ws.Range("D46").Value = "'43.88"
$ws.Range("E46").Value = "  -0.92%  "

$ws.Range("E47").Value = "  -4.47%  "

$ws.Range("D48").Value = "'166.19"
$ws.Range("E48").Value = "  -8.67%  "

$ws.Range("D49").Value = "'0.725"
$ws.Range("E49").Value = "  -1.11%  "

$ws.Range("D50").Value = "'24.67"
$ws.Range("E50").Value = "  -1.72%  "

$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").Value = "'0.600"
$ws.Range("E51").Value = "  -5.69%  "

